$d = $word.ActiveDocument

# 1. "Version 1." -> "Version 2" (the trailing "." is re-added afterwards
#    as its own run, positioned after the _GoBack bookmark, matching the
#    original layout). Only "1." is touched so the spellStart/spellEnd
#    proofErr markers around "Version" stay put.
[void]$d.Content.Find.Execute("1.", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "2", 2)

# 2. Split the run "Version" into "Versi" | "on" at the boundary after
#    the 5th character, without leaving any residual formatting (rPr).
#    Adding then deleting a bookmark at that exact point forces the run
#    to break there cleanly.
$splitPoint = $d.Range(5, 5)
[void]$d.Bookmarks.Add("__tmp_run_split__", $splitPoint)
$d.Bookmarks("__tmp_run_split__").Delete()

# 3. Re-add the final "." as its own run, positioned after the _GoBack
#    bookmark (matching the original document layout).
$goBack = $d.Bookmarks("_GoBack")
[void]$goBack.Range.InsertAfter(".")
